# Update "docs/Datenbankstruktur.xlsx" for player_matches table:
# add oks_atk / ods_atk columns (V10/W10), and clear out the
# now-unused reference table in column C (rows 34-54), leaving
# only the two bolded separator cells (C34, C45) behind, blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New header cells on the PlayerMatch row.
$ws.Range("V10").Value = "oks_atk"
$ws.Range("W10").Value = "ods_atk"

# Clear the old column-C reference list entirely (rows 34-54), then
# restore the bold style on the two rows that still remain (34, 45),
# leaving their contents blank.
$ws.Range("C34:C54").ClearContents()
$ws.Range("C34").Font.Bold = $true
$ws.Range("C45").Font.Bold = $true

# Update the view: scroll so column K is left-most and select W10.
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("W10").Select() | Out-Null
